$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1: 土地 (land) — add data row 2, index 14
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("A2").Value = 14
$ws1.Range("B1").Copy()
$ws1.Range("A2").PasteSpecial(-4122)

$ws1.Range("B2").Value = "臺北市大安區龍泉段一小段02930000地號"
$ws1.Range("C2").Value = 365
$ws1.Range("D2").Value = "100000分之16216"
$ws1.Range("E2").Value = "高金素梅"
$ws1.Range("F2").Value = "92年12月25日"
$ws1.Range("G2").Value = "033貝賣"
$ws1.Range("H2").Value = "25000000(土地建物與車位合併價）"
$ws1.Range("I2").Value = "land"
$ws1.Range("J2").Value = "normal"
$ws1.Range("K2").Value = "'2012-04-30"
$ws1.Range("L2").Value = "高金素梅"
$ws1.Range("M2").Value = 926
$ws1.Range("N2").Value = "tmpb18e1"
$ws1.Range("O2").Value = 14
$ws1.Range("P2").Value = 0.16216
$ws1.Range("Q2").Value = 59.1884

# ---------------------------------------------------------------
# Sheet 2: 建物 (building) — row 1 held data-shaped values (a quirk
# of the source generator); it must become a proper header row like
# sheet 1's, and the real data moves to row 2 (index 19).
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("B1").Value = "name"
$ws2.Range("C1").Value = "area"
$ws2.Range("D1").Value = "share_portion"
$ws2.Range("E1").Value = "owner"
$ws2.Range("F1").Value = "register_date"
$ws2.Range("G1").Value = "register_reason"
$ws2.Range("H1").Value = "acquire_value"
$ws2.Range("I1").Value = "property_category"
$ws2.Range("J1").Value = "category"
$ws2.Range("K1").Value = "date"
$ws2.Range("L1").Value = "legislator_name"
$ws2.Range("M1").Value = "legislator_id"
$ws2.Range("N1").Value = "source_file"
$ws2.Range("O1").Value = "index"
$ws2.Range("P1").Value = "portion"
$ws2.Range("Q1").Value = "total"

# The I1:Q1 cells are brand new — give them the same bold/bordered
# header style already used by B1:H1 (style index 1).
$ws2.Range("B1").Copy()
$ws2.Range("I1:Q1").PasteSpecial(-4122)

$ws2.Range("A2").Value = 19
$ws1.Range("B1").Copy()
$ws2.Range("A2").PasteSpecial(-4122)

$ws2.Range("B2").Value = "臺北市大安區龍泉段一小段05819000建號"
$ws2.Range("C2").Value = 148.31
$ws2.Range("D2").Value = "全部"
$ws2.Range("E2").Value = "高金素梅"
$ws2.Range("F2").Value = "92年12月25日"
$ws2.Range("G2").Value = "貝賣"
$ws2.Range("H2").Value = "25000000(土地建物與車位合併價）"
$ws2.Range("I2").Value = "land"
$ws2.Range("J2").Value = "normal"
$ws2.Range("K2").Value = "'2012-04-30"
$ws2.Range("L2").Value = "高金素梅"
$ws2.Range("M2").Value = 926
$ws2.Range("N2").Value = "tmpb18e1"
$ws2.Range("O2").Value = 19
$ws2.Range("P2").Value = 1
$ws2.Range("Q2").Value = 148.31

# ---------------------------------------------------------------
# Sheet 3: 存款 (deposit) — add data row 2, index 45
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("A2").Value = 45
$ws1.Range("B1").Copy()
$ws3.Range("A2").PasteSpecial(-4122)

$ws3.Range("B2").Value = "臺灣銀行群賢分行"
$ws3.Range("C2").Value = "活期存款"
$ws3.Range("D2").Value = "新臺幣"
$ws3.Range("E2").Value = "高金素梅"
$ws3.Range("F2").Value = 366738

# ---------------------------------------------------------------
# Sheet 4: 債務 (debt) — regenerated source data:
#   * row 1 (the pseudo-header) now mirrors the new first record
#   * a new row 2 (index 85) is inserted
#   * the former row 2 (index 86) becomes row 3, with a couple of
#     values swapped in for the new second record
# ---------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

# Move the old row-2 record (index 86) down to row 3 first, copying
# its style so A3 keeps the bold/bordered "index" look.
$ws4.Range("A2").Copy()
$ws4.Range("A3").PasteSpecial(-4122)

$ws4.Range("A3").Value = 86
$ws4.Range("B3").Value = "現金"
$ws4.Range("C3").Value = "局金素梅"
$ws4.Range("D3").Value = "石旭松新北市泰山區明志路"
$ws4.Range("E3").Value = "'4000000"
$ws4.Range("F3").Value = "96年02月06日"
$ws4.Range("G3").Value = "借款"

# New row 2 (index 85).
$ws4.Range("A2").Value = 85
$ws4.Range("B2").Value = "現金"
$ws4.Range("C2").Value = "高金素梅"
$ws4.Range("D2").Value = "陳麗卿新北市泰山區明志路"
$ws4.Range("E2").Value = "'6000000"
$ws4.Range("F2").Value = "96年02月06日"
$ws4.Range("G2").Value = "借款"

# Row 1 (pseudo-header) now mirrors the new row-2 record.
$ws4.Range("B1").Value = "現金"
$ws4.Range("C1").Value = "高金素梅"
$ws4.Range("D1").Value = "陳麗卿新北市泰山區明志路"
$ws4.Range("E1").Value = "'6000000"
$ws4.Range("F1").Value = "96年02月06日"
$ws4.Range("G1").Value = "借款"
